# Cap nhat cong viec den 27/10
# Refresh the Task/Result/Assignment columns for weeks 2-8, re-flow the
# "Ket qua" column (E) to wrap like the rest of the table, and extend
# wrapped/centered formatting into the newly-used part of "Phan cong" (F).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Formats: clone existing styles onto the ranges that need them,  ---
# --- instead of toggling individual alignment flags one at a time.   ---

# E1 ("Ket qua" header) -> wrap-only style already used by D1.
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null

# E2:E12 -> center+wrap style already used throughout column D.
$ws.Range("D2").Copy() | Out-Null
$ws.Range("E2:E12").PasteSpecial(-4122) | Out-Null

# F2 -> center-only style already used in column A.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("F2").PasteSpecial(-4122) | Out-Null

# F3:F5 -> center+wrap style (same one used in column D).
$ws.Range("D2").Copy() | Out-Null
$ws.Range("F3:F5").PasteSpecial(-4122) | Out-Null

# F6:F9 -> center-only style for the newly-used blank assignment cells.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("F6:F9").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# Widen column E for the longer wrapped results (replaces the old
# "best fit" auto width).
$ws.Columns.Item(5).ColumnWidth = 44.7213541666667

# --- Content updates ---

# Row 3: usecase/prototype finished.
$ws.Range("E3").Value = "hoàn thành usecase, sequence, class diagram, lược đồ csdl, prototype"

# Row 4: new task + result text, assignment note tweak.
$ws.Range("D4").Value = "thiết kế UI cho app, web, viết api"
$ws.Range("E4").Value = "Hoàn thành xong giao diện web`nApp: hoàn thành giao diện main, login, person"
$ws.Range("F4").Value = "Thọ: Giao diện phần app.`nCông: Giao diện phần web, viết API."

# Row 5: same task text as row 4 now; assignment note carries down too.
$ws.Range("D5").Value = "thiết kế UI cho app, web, viết api"
$ws.Range("F5").Value = "Thọ: Giao diện phần app.`nCông: Giao diện phần web, viết API."

# Rows 6-7: task text back to "Code".
$ws.Range("D6").Value = "Code"
$ws.Range("D7").Value = "Code"

# --- Row heights recomputed by Excel for the new wrapped content ---
$ws.Rows.Item(4).RowHeight = 54
$ws.Rows.Item(5).RowHeight = 36

# --- View: zoom to 85% and move the active selection to E12 ---
$excel.ActiveWindow.Zoom = 85
$ws.Range("E12").Select() | Out-Null
